$wb = $excel.ActiveWorkbook

# Layer0 sheet updates
$ws0 = $wb.Worksheets.Item("Layer0")
$ws0.Range("B2").Value = -0.6207974326880873
$ws0.Range("C2").Value = -0.2936463303025039
$ws0.Range("B3").Value = -0.008752254176590913
$ws0.Range("C3").Value = 0.9762913825547721
$ws0.Range("B4").Value = 0.8681185898850409
$ws0.Range("C4").Value = 0.06044466525304374

# Layer1 sheet updates
$ws1 = $wb.Worksheets.Item("Layer1")
$ws1.Range("B2").Value = -1.754117842938022
$ws1.Range("C2").Value = -0.3558846413607962
$ws1.Range("B3").Value = 0.7087473255926309
$ws1.Range("C3").Value = 0.1880960260198566
$ws1.Range("B4").Value = 0.8867000596440289
$ws1.Range("C4").Value = 0.4744587374238869
